$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 84.9697536618426
$ws.Range("B3").Value = 88.50669163758266
$ws.Range("B4").Value = 91.31851568575865
$ws.Range("H5").Value = 95.88612394129422
$ws.Range("H6").Value = 95.88045347654935
$ws.Range("H7").Value = 95.84900992240568
$ws.Range("C8").Value = 98.55771089523063
$ws.Range("C9").Value = 97.24678051129698
$ws.Range("C10").Value = 98.1894929294416
$ws.Range("D11").Value = 99.30358666383916
$ws.Range("D12").Value = 99.1991717128687
$ws.Range("D13").Value = 99.2520149104038
$ws.Range("E14").Value = 98.71533270856034
$ws.Range("E15").Value = 98.77466933319627
$ws.Range("E16").Value = 98.73644991363936
$ws.Range("F17").Value = 98.20976544239495
$ws.Range("F18").Value = 98.27483948254691
$ws.Range("F19").Value = 98.19319468011865
$ws.Range("G20").Value = 97.24264997032257
$ws.Range("G21").Value = 97.33136547552698
$ws.Range("G22").Value = 97.27608855599398
$ws.Range("B23").Value = 90.10870807531742
$ws.Range("B24").Value = 93.15718638210201
$ws.Range("H25").Value = 95.9143479783671
$ws.Range("H26").Value = 95.86307539611511
$ws.Range("C27").Value = 97.81429708426059
$ws.Range("C28").Value = 98.02049535525029
$ws.Range("D29").Value = 99.26492918236333
$ws.Range("D30").Value = 99.25353462635834
$ws.Range("E31").Value = 98.76337678055769
$ws.Range("E32").Value = 98.72366983585701
$ws.Range("F33").Value = 98.23910002563062
$ws.Range("F34").Value = 98.30770226925009
$ws.Range("G35").Value = 97.20979548381735
$ws.Range("G36").Value = 97.22904729439165
$ws.Range("B37").Value = 91.63091943761026
$ws.Range("B38").Value = 94.26170851986157
$ws.Range("H39").Value = 95.86309660364532
$ws.Range("H40").Value = 95.8095743975328
$ws.Range("C41").Value = 98.56853365522628
$ws.Range("C42").Value = 98.26427319361977
$ws.Range("D43").Value = 99.27167091385614
$ws.Range("D44").Value = 99.21734387641084
$ws.Range("E45").Value = 98.71442669635556
$ws.Range("E46").Value = 98.65434657647187
$ws.Range("F47").Value = 98.273671367419
$ws.Range("F48").Value = 98.14526655144304
$ws.Range("G49").Value = 97.22988135560763
$ws.Range("G50").Value = 97.33869547284199
